$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Change the Results column (D) for rows 2-4 (TestCase_E1, E2, E3)
# from "PASS" to "SKIP", matching the updated run-mode results.
$ws.Range("D2:D4").Value = "SKIP"
